$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (price) cells to text before writing, so values like
# "595.38" or "1.00" are not auto-converted to numbers (they were
# originally stored as inlineStr text). ClearFormats() afterwards drops
# the temporary "@" number-format so no stray style survives on the cell.
# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.711.19'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.64%  '
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.777.94'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +1.32%  '
# Row 4
$ws.Range("E4").Value = '  +0.00%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.38'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.65%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.90'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.45%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.764.60'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +1.01%  '
# Row 8
$ws.Range("E8").Value = '  +0.05%  '
# Row 9
$ws.Range("E9").Value = '  +0.61%  '
# Row 10
$ws.Range("E10").Value = '  -0.10%  '
# Row 11
$ws.Range("E11").Value = '  -2.15%  '
# Row 12
$ws.Range("E12").Value = '  +0.22%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000253'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -2.00%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.03'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.17%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.413.03'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.21%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.777.65'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.13%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.53'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +3.82%  '
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.642.66'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.83%  '
# Row 19
$ws.Range("E19").Value = '  +0.41%  '
# Row 20
$ws.Range("E20").Value = '  +0.04%  '
# Row 21
$ws.Range("E21").Value = '  -5.75%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '459.49'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.16%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.697'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.32%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000155'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +4.63%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.39'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.38%  '
# Row 26
$ws.Range("E26").Value = '  +1.36%  '
# Row 27
$ws.Range("E27").Value = '  -2.20%  '
# Row 28
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.02'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.06%  '
# Row 29
$ws.Range("B29").Value = 'Dai'
$ws.Range("C29").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.10%  '
# Row 30
$ws.Range("E30").Value = '  +0.57%  '
# Row 31
$ws.Range("E31").Value = '  +3.88%  '
# Row 32
$ws.Range("E32").Value = '  -0.68%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '29.54'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.40%  '
# Row 34
$ws.Range("E34").Value = '  -0.14%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.08'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.75%  '
# Row 36
$ws.Range("E36").Value = '  -0.02%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.38'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.09%  '
# Row 38
$ws.Range("E38").Value = '  +0.01%  '
# Row 39
$ws.Range("E39").Value = '  +0.06%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.76'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.14%  '
# Row 41
$ws.Range("E41").Value = '  -0.12%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '45.52'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +4.02%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '48.06'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +3.39%  '
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.298'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.01%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '149.86'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +3.87%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.30'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.52%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '393.72'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.54%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '26.64'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +6.31%  '
# Row 50
$ws.Range("E50").Value = '  -4.70%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.719.32'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.79%  '
